# Update the "time_taken" column (F) on the "data" sheet with refreshed
# query timestamps, then add a new "metadata" worksheet (after "data")
# describing the panel query itself.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:33:54.449674",
    "2021-10-05 14:33:54.449681",
    "2021-10-05 14:33:54.449684",
    "2021-10-05 14:33:54.449687",
    "2021-10-05 14:33:54.449690",
    "2021-10-05 14:33:54.449692",
    "2021-10-05 14:33:54.449695",
    "2021-10-05 14:33:54.449697",
    "2021-10-05 14:33:54.449700",
    "2021-10-05 14:33:54.449703",
    "2021-10-05 14:33:54.449705",
    "2021-10-05 14:33:54.449708",
    "2021-10-05 14:33:54.449710",
    "2021-10-05 14:33:54.449712",
    "2021-10-05 14:33:54.449715",
    "2021-10-05 14:33:54.449717",
    "2021-10-05 14:33:54.449720",
    "2021-10-05 14:33:54.449723",
    "2021-10-05 14:33:54.449725",
    "2021-10-05 14:33:54.449728",
    "2021-10-05 14:33:54.449730",
    "2021-10-05 14:33:54.449733",
    "2021-10-05 14:33:54.449735",
    "2021-10-05 14:33:54.449738",
    "2021-10-05 14:33:54.449741",
    "2021-10-05 14:33:54.449743",
    "2021-10-05 14:33:54.449746",
    "2021-10-05 14:33:54.449748",
    "2021-10-05 14:33:54.449750"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Add the new "metadata" sheet right after "data".
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Glycogen Storage Diseases"
$meta.Cells.Item(2, 3).Value = 106
# Keep "1.1" as text (matching the source panel version string) instead of
# letting Excel coerce it to the numeric value 1.1.
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.1"
$meta.Cells.Item(2, 5).Value = "2021-04-08T10:48:29.657674Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:33:54.446094"
$meta.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/106/?format=json"

# Match the header / index-column styling used on the "data" sheet
# (bold, bordered, centered alignment).
$headerRange = $meta.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$indexCell = $meta.Cells.Item(2, 1)
$indexCell.Font.Bold = $true
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160
$indexCell.Borders.LineStyle = 1
